# businessNeed.xlsx — add the "Audio Recording" business-need row (sourced
# from the refreshed SharePoint "Business Need" list) and pick up a handful
# of text corrections/additions on existing rows that came in with the same
# refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Insert a new data row right under the header/filter row (row 2), --
# --- shifting "Brainstorming" and everything below it down by one.    --
$ws.Rows(3).Insert()

# The table (ListObject) + its AutoFilter need to grow to include the new
# row; Insert() alone doesn't resize the table definition.
$lo.Resize($ws.Range("A1:F17"))

# Keep the hidden workbook-scoped defined name (used by the query table)
# in sync with the new table extent.
$n = $wb.Names.Item(1)
$n.RefersTo = "=owssvr!`$A`$1:`$F`$17"

# --- Populate the new row 3: "Audio Recording" ---
$ws.Cells.Item(3, 1).Value = "Audio Recording"
$ws.Cells.Item(3, 2).Value = "_Starting Over;#24;#Adobe Audition;#54"
$ws.Cells.Item(3, 3).Value = "_Starting Over;#22;#Audio;#61"
$ws.Cells.Item(3, 5).Value = "Item"
$ws.Cells.Item(3, 6).Value = "teams/kmqa/Lists/Business Need"

# --- Text fixes/additions on existing rows (row numbers below are the ---
# --- NEW (post-insert) row numbers).                                   --

# Communication (row 6): "New Letters" -> "Newsletter"
$ws.Cells.Item(6, 3).Value = "_Starting Over;#22;#InfoGraphic;#47;#Posters;#48;#Newsletter;#49;#Memes and Vemes;#50"

# Learning (row 12): "Adobe Acrobat 9 Pro" -> "Adobe Acrobat Pro", plus
# new trailing software/deliverable options.
$ws.Cells.Item(12, 2).Value = "_Starting Over;#24;#Adobe Acrobat Pro;#34;#Adobe eLearning;#35;#Presenter Media;#36;#Prezi;#37;#Prespectore;#38;#Microsoft PowerPoint;#39;#Microsoft Word;#50;#Adobe Lifecycle;#51;#Designer ES2;#52;#Microsoft Publisher;#53"
$ws.Cells.Item(12, 3).Value = "_Starting Over;#22;#Curriculum;#34;#eBook;#35;#Job Aid;#36;#Online Module;#37;#Presentations;#38;#Training Manual;#39;#Information Sheet;#60"

# Research (row 14): new trailing deliverable options.
$ws.Cells.Item(14, 3).Value = "_Starting Over;#22;#Qualitative and Quantitative Analysis;#42;#Knowledge Repositories;#43;#Literature Review;#44;#Environmental Scan;#59;#Literature Search;#62;#Systematic Review;#63;#Evidence Summary;#64"

# Videography (row 17): new trailing software option.
$ws.Cells.Item(17, 2).Value = "_Starting Over;#24;#Soney Vegas;#47;#Go Animate;#55"
